# Add 2022-Q3 data
# 1) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 at the top
#    of the data (row 2), push existing quarters down, and renumber the
#    leading index column.
$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")

# Insert a new blank row right below the header row.
$summary.Rows.Item(2).Insert()

# Match the formatting used by the rest of column A (bold, thin box border,
# centered / top aligned) since Insert() does not carry it onto the new cell;
# copy it wholesale from the cell directly below so it resolves to the exact
# same style record instead of constructing a new (slightly different) one.
$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# Fill the new row with the 2022-Q3 totals.
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 1
$summary.Cells.Item(2, 4).Value = 0.06

# Renumber the index column for the rows that shifted down (they kept their
# old 0-based index values, so bump each by one).
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(9, 1).Value = 7

# 2) Insert a brand new "2022-Q3" sheet before the existing "2022-Q2" sheet
#    holding the fund holdings detail for the new quarter.
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

# Seed the header row + the index-column style by copying them from the
# existing "2022-Q2" sheet, which already uses the workbook's shared header
# style - this reuses that exact style record instead of creating a new one.
$q2.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q2.Cells.Item(2, 1).Copy()
$q3.Cells.Item(2, 1).PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$q3.Cells.Item(2, 1).Value = 0

# These columns hold text-formatted numbers in every other quarter sheet, so
# force text formatting before assigning the values.
$q3.Cells.Item(2, 2).NumberFormat = "@"
$q3.Cells.Item(2, 2).Value = "515450"
$q3.Cells.Item(2, 3).NumberFormat = "@"
$q3.Cells.Item(2, 3).Value = "南方标普中国A股大盘红利低波50ETF"
$q3.Cells.Item(2, 4).NumberFormat = "@"
$q3.Cells.Item(2, 4).Value = "2.17"
$q3.Cells.Item(2, 5).NumberFormat = "@"
$q3.Cells.Item(2, 5).Value = "99.66"
$q3.Cells.Item(2, 6).NumberFormat = "@"
$q3.Cells.Item(2, 6).Value = "2.65"
$q3.Cells.Item(2, 7).NumberFormat = "@"
$q3.Cells.Item(2, 7).Value = "0.0575"
$q3.Cells.Item(2, 8).Value = 8
